$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$header = $ws.Range("A1:F1")
$header.Style = "Normal"

for ($r = 2; $r -le 19; $r++) {
    $c = $ws.Cells.Item($r, 3)
    $c.Style = "Normal"
}

$ws.Range("A20").Value = "http://localhost//shopping/admin/productimages/21/SamsungTVFrontView.jpg"
$ws.Range("B20").Value = "http://localhost//shopping/product-details.php?pid=21"
$ws.Range("D20").Value = "Rs.139900"
$ws.Range("E20").Value = "Rs.0"
$ws.Range("F20").Value = "Add to Cart"

$ws.Range("A21").Value = "http://localhost//shopping/admin/productimages/22/SamsungTVFrontView.jpg"
$ws.Range("B21").Value = "http://localhost//shopping/product-details.php?pid=22"
$ws.Range("D21").Value = "Rs.139900"
$ws.Range("E21").Value = "Rs.0"
$ws.Range("F21").Value = "Add to Cart"

$c20 = $ws.Cells.Item(20, 3)
$c20.Style = "Normal"

$c21 = $ws.Cells.Item(21, 3)
$c21.Style = "Normal"
